# Generate Report for Archive
# - Update the "Status" value from "Ready for handoff" to "In Translation"
#   everywhere it appears (Overview!E2/F2, zh-cn!C2, de-de!C2 all share the
#   same string).
# - Shrink the "Status" column(s) to fit the new, shorter text
#   (Overview columns E & F, and column C on the zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
